$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.177.00'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '1.644.99'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = "'217.47"
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("D6").Value = "'0.510"
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").Value = "'19.96"
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("D12").Value = '1.873.73'
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("D13").Value = '1.644.07'
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("E14").Value = '  -2.44%  '
$ws.Range("D15").Value = "'0.540"
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = "'67.44"
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '27.135.33'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").Value = "'218.52"
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = "'6.85"
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("D23").Value = "'2.51"
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("D24").Value = "'9.19"
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("D25").Value = "'147.54"
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").Value = "'7.43"
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("D29").Value = "'15.75"
$ws.Range("E29").Value = '  -2.04%  '
$ws.Range("D30").Value = "'0.0504"
$ws.Range("E30").Value = '  -2.18%  '
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").Value = "'3.38"
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("D35").Value = '1.268.06'
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").Value = "'0.0178"
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("E38").Value = '  +0.24%  '
$ws.Range("D39").Value = "'0.840"
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("D41").Value = "'0.808"
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("E42").Value = '  +4.30%  '
$ws.Range("D43").Value = "'5.41"
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").Value = '1.784.98'
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").Value = "'91.75"
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("E48").Value = '  +14.46%  '
$ws.Range("D49").Value = "'0.0512"
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("D50").Value = "'7.65"
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").Value = "'0.0975"
$ws.Range("E51").Value = '  -1.04%  '
